$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# --- e009 "Ammo Loading Limits" (row 10, column B) ---------------------
# Text reworded: "added in a later step" now ends the sentence directly
# instead of pointing at the e009b button, and "--AP:" gains a space to
# become "-- AP:" (matching the existing "-- HE:" line).
$e009 = @'
<Bold>e009 Ammo Loading Limits</Bold> <InlineUIContainer><Button Content='r16.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
See 
<InlineUIContainer><Button Content='r16.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
for ammo types. See 
<InlineUIContainer><Button Content='r16.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
for loading ammo. The Tank Card limits the number of normal main gun ammo allowed to AMMO_NORMAL_LOAD. Extra ammo is added in a later step after assigning normal load.
<LineBreak/><LineBreak/>
<Bold>-- AP:</Bold> Unlimited<LineBreak/>
<Bold>-- HE:</Bold> Unlimited
'@
$ws.Range("B10").Value = $e009.TrimEnd()

# e009's row got one line shorter -> row height shrinks from 165 to 150.
$ws.Rows.Item(10).RowHeight = 150

# --- e028 "Enter Adjacent Area" (row 29, column B) ----------------------
# Swap the placeholder "Combat" image for the real "Sherman1" artwork
# (bigger + re-indented) that illustrates this rule section.
$e028 = @'
<Bold>e028 Enter Adjacent Area</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Click on one of the adjacent highlighted areas. Artillery Support or Air Strike Counters are moved to the battle board as a reminder. 
<LineBreak/><LineBreak/>
                        <InlineUIContainer><Image Name='Sherman1' Height='200' Width='325'></Image></InlineUIContainer>
'@
$ws.Range("B29").Value = $e028.TrimEnd()

# --- Scroll / selection bookkeeping (matches the saved view) -----------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("B30").Select()
